# Adiciona a função de comparar o aluno com a média do curso
# Acrescenta duas novas linhas de log de acesso (linhas 99 e 100) à planilha.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A coluna "data" (B) contém valores como "2025-06-03" que o Excel tentaria
# interpretar como datas. Forçamos formato de texto antes de escrever os
# valores para mante-los como texto simples, depois removemos a formatação
# para que a célula fique sem estilo explicito (igual as demais linhas).
$dataRange = $ws.Range("B99:B100")
$dataRange.NumberFormat = "@"

$ws.Range("A99").Value = "coord123"
$ws.Range("B99").Value = "2025-06-03"
$ws.Range("C99").Value = "01:38:39"

$ws.Range("A100").Value = "coord123"
$ws.Range("B100").Value = "2025-06-03"
$ws.Range("C100").Value = "02:03:03"

$dataRange.ClearFormats()
